# Convert the "isMissing" column (column C) on the "Categories" sheet
# from the text string "FALSE" to the real boolean value FALSE for every
# data row (rows 2-301).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Categories")
$ws.Range("C2:C301").Value = $false
